# Ivy CHUANG_CampB_timetable.xlsx edit:
# Shift the column-G event labels (and their merged blocks) down by one row,
# and swap the custom row heights between the 26/27 and 44/45 row pairs so
# the taller row follows the shifted label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unmerge all existing column-G merged blocks so the new merges below
#        don't collide with (and silently absorb) the old boundaries. ---
$ws.Range("G7:G13").UnMerge()
$ws.Range("G14:G19").UnMerge()
$ws.Range("G20:G25").UnMerge()
$ws.Range("G26:G28").UnMerge()
$ws.Range("G29:G31").UnMerge()
$ws.Range("G32:G35").UnMerge()
$ws.Range("G36:G40").UnMerge()
$ws.Range("G41:G43").UnMerge()
$ws.Range("G44:G50").UnMerge()

# --- 2. Move each label text down to the next row in column G, clearing the
#        old cell. Text is written out literally (rather than copied via
#        .Value2) so the original "\r\n" line breaks are preserved exactly. ---
$ws.Range("G15").Value = "Lunch"
$ws.Range("G14").Value = ""

$ws.Range("G21").Value = "Lina Summer Camp of Music Students & Friends Concert"
$ws.Range("G20").Value = ""

$ws.Range("G27").Value = "After concert refreshment `r`n(Maritime Museum)"
$ws.Range("G26").Value = ""

$ws.Range("G30").Value = "Rehearsal for Faculty Concert"
$ws.Range("G29").Value = ""

$ws.Range("G33").Value = "Break"
$ws.Range("G32").Value = ""

$ws.Range("G37").Value = "Lina Summer Camp of Music Faculty Concert"
$ws.Range("G36").Value = ""

$ws.Range("G42").Value = "After concert refreshment"
$ws.Range("G41").Value = ""

$ws.Range("G45").Value = "After Concert Dinner `r`n(Pheasant-Jasmine Room, Mandarin Oriental)"
$ws.Range("G44").Value = ""

# --- 3. Swap the custom row heights so the taller (30pt) row now lands on
#        the row that holds the moved label. ---
$ws.Rows.Item(26).RowHeight = 15
$ws.Rows.Item(27).RowHeight = 30

$ws.Rows.Item(44).RowHeight = 15
$ws.Rows.Item(45).RowHeight = 30

# --- 4. Re-create the column-G merges, each shifted down by one row. Re-apply
#        a uniform thin box border afterwards -- Merge() otherwise splits the
#        border definition per interior row, whereas the sheet's convention
#        (see the untouched G3:G6 merge) is a single thin border style
#        reused verbatim on every cell of a merged block. ---
$ws.Range("G7:G14").Merge()
$ws.Range("G7:G14").Borders.LineStyle = 1

$ws.Range("G15:G20").Merge()
$ws.Range("G15:G20").Borders.LineStyle = 1

$ws.Range("G21:G26").Merge()
$ws.Range("G21:G26").Borders.LineStyle = 1

$ws.Range("G27:G29").Merge()
$ws.Range("G27:G29").Borders.LineStyle = 1

$ws.Range("G30:G32").Merge()
$ws.Range("G30:G32").Borders.LineStyle = 1

$ws.Range("G33:G36").Merge()
$ws.Range("G33:G36").Borders.LineStyle = 1

$ws.Range("G37:G41").Merge()
$ws.Range("G37:G41").Borders.LineStyle = 1

$ws.Range("G42:G44").Merge()
$ws.Range("G42:G44").Borders.LineStyle = 1

$ws.Range("G45:G50").Merge()
$ws.Range("G45:G50").Borders.LineStyle = 1
